$d = $word.ActiveDocument

$replacements = @(
    @("982×3=2946", "155×9=1395"),
    @("133×6=798", "875×6=5250"),
    @("703×5=3515", "191×4=764"),
    @("466×9=4194", "271×5=1355"),
    @("793×4=3172", "513×9=4617"),
    @("775×7=5425", "569×3=1707"),
    @("202×4=808", "627×6=3762"),
    @("589×4=2356", "564×5=2820"),
    @("791×7=5537", "224×8=1792"),
    @("130×4=520", "658×4=2632"),
    @("630×7=4410", "892×4=3568"),
    @("483×3=1449", "914×5=4570"),
    @("480×4=1920", "585×9=5265"),
    @("571×6=3426", "722×6=4332"),
    @("656×8=5248", "347×9=3123"),
    @("268×6=1608", "441×9=3969"),
    @("624×8=4992", "275×3=825"),
    @("509×9=4581", "915×7=6405"),
    @("741×2=1482", "737×5=3685"),
    @("677×4=2708", "684×5=3420"),
    @("617×3=1851", "919×9=8271"),
    @("267×3=801", "691×6=4146"),
    @("185×7=1295", "507×5=2535"),
    @("970×6=5820", "284×5=1420"),
    @("985×6=5910", "175×4=700")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
